$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R with year 2021 header and its value
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 102.20441221981518

# Copy style of Q4 (year header) to R4
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Apply custom number format "0.0" with right/center alignment + border to R5,
# matching the existing data-row style but with its own number format.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R5").NumberFormat = "0.0"

# update selection to match diff
$ws.Range("S9").Select() | Out-Null
